# RoStatProcessing.xlsx edit:
# Damage indicator pooling and some cast effect stuff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# StatDef sheet: assorted monster stat tweaks
# ---------------------------------------------------------------------
$statDef = $wb.Worksheets.Item("StatDef")

$statDef.Range("D19").Value = 19
$statDef.Range("D52").Value = 10
$statDef.Range("AC81").Value = "Golem"
$statDef.Range("AH139").Value = 1
$statDef.Range("AC173").Value = "Small"
$statDef.Range("AC174").Value = "Brute"
$statDef.Range("AC175").Value = "Brute,Elite"
$statDef.Range("AC179").Value = "Demon"
$statDef.Range("AC180").Value = "Undead"
$statDef.Range("AH212").Value = 1
$statDef.Range("AH214").Value = 1
$statDef.Range("AH217").Value = 1
$statDef.Range("AH218").Value = 1
$statDef.Range("AH222").Value = 1
$statDef.Range("H250").Value = 140
$statDef.Range("O250").Value = 130
$statDef.Range("Q250").Value = 110
$statDef.Range("AH264").Value = 1
$statDef.Range("AC265").Value = "Strong"
$statDef.Range("AH269").Value = -1
$statDef.Range("AC270").Value = "Brute"
$statDef.Range("AH283").Value = 1
$statDef.Range("AC285").Value = "Normal,Ranged"
$statDef.Range("AH285").Value = -1

# ---------------------------------------------------------------------
# ClassDef sheet: insert a new "Small" class row before "Ranged"
# ---------------------------------------------------------------------
$classDef = $wb.Worksheets.Item("ClassDef")

$classDef.Rows.Item(6).Insert()
$classDef.Range("A6").Value = "Small"
$classDef.Range("B6:M6").Value = 70
$classDef.Range("N6").Value = 1

# ---------------------------------------------------------------------
# Window / selection state: ClassDef becomes the active sheet
# ---------------------------------------------------------------------
$statDef.Range("C154").Select()
$statDef.Range("AC173").Select()

$classDef.Activate()
$classDef.Range("B6").Select()
